$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values -----------------------------------------------------
# A1: link text, A2: page title text, A3: cleared out (was "ghi")
$ws.Range("A1").Value = "https://www.tsx.com"
$ws.Range("A2").Value = "TMX TSX | TSXV - Toronto Stock Exchange and TSX Venture Exchange"
$ws.Range("A3").Value = ""

# --- Formatting --------------------------------------------------------
# A1 and A2 now wrap their text (matches the new shared cellXf 1)
$ws.Range("A1").WrapText = $true
$ws.Range("A2").WrapText = $true

# A3 gets the "Normal" cell style re-applied explicitly (so a distinct,
# protection-flagged cellXf is written instead of the cell being dropped
# from sheetData entirely).
$ws.Range("A3").Style = "Normal"

# Column A is widened to fit the long title text
$ws.Columns.Item(1).ColumnWidth = 69.75

# --- Selection -----------------------------------------------------------
$ws.Range("A1").Select()

$wb.Save()
